$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.037.47"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "2.090.00"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "228.85"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D7").Value = "60.97"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  +2.74%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "2.399.58"
$ws.Range("E12").Value = "  +2.93%  "
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("E14").Value = "  +3.56%  "
$ws.Range("E15").Value = "  +6.35%  "
$ws.Range("D16").Value = "0.774"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").Value = "2.098.67"
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("D18").Value = "37.983.15"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("D20").Value = "69.99"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "0.0₃0838"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").Value = "224.02"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  +2.45%  "
$ws.Range("D26").Value = "169.55"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").Value = "9.38"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("E28").Value = "  +3.33%  "
$ws.Range("D29").Value = "18.95"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +4.24%  "
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").Value = "2.38"
$ws.Range("E32").Value = "  +10.70%  "
$ws.Range("D33").Value = "4.67"
$ws.Range("E33").Value = "  +3.57%  "
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").Value = "0.0605"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").Value = "2.41"
$ws.Range("E36").Value = "  +5.61%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "3.55"
$ws.Range("E38").Value = "  +8.92%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "18.06"
$ws.Range("E40").Value = "  +4.80%  "
$ws.Range("D41").Value = "1.544.23"
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").Value = "99.92"
$ws.Range("E42").Value = "  +3.92%  "
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").Value = "4.17"
$ws.Range("E46").Value = "  +5.51%  "
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").Value = "2.286.60"
$ws.Range("E51").Value = "  +2.91%  "
